$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format so numeric-looking strings (e.g. "213.00", "1.01")
# are preserved exactly instead of being auto-converted to numbers.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.731.53'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.633.04'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.54%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '219.35'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.496'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -1.35%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.56%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -1.31%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -1.38%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '18.96'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0844'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.859.24'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.86%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.628.61'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -1.41%  '
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -2.07%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.521'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -1.32%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.23'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -1.03%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.706.54'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '213.00'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.84%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.23'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.30'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -4.66%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.02'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -3.67%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '147.56'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.53%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.34%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.59'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.80%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0502'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -2.79%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.20'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.44%  '
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -1.07%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.252.18'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -1.85%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0175'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -2.89%  '
$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.805'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -2.77%  '
$ws.Range('B40').NumberFormat = "@"
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').NumberFormat = "@"
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.01'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -1.60%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.26'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.76%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.770.18'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.11'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -5.68%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '91.80'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '59.57'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -2.17%  '
$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0515'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.76%  '
$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'USDD'
$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.01'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0955'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.99%  '
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -0.81%  '
